# Non-Oncology Excel reports comparison
# Update the "ExcelReport" expected-filename strings so the hyphen after
# "NewImportLogic_1" is no longer padded with spaces on either side
# (" - " -> "-"), matching the WordReport/StandardExcelReport siblings
# that were already in that format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value  = "ExcelReport-NewImportLogic_1-Test_Automation_1-Clinical-"
$ws.Range("G6").Value  = "ExcelReport-NewImportLogic_1-Test_Automation_1-Economic-"
$ws.Range("G9").Value  = "ExcelReport-NewImportLogic_1-Test_Automation_1-Quality of Life-"
$ws.Range("G12").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Real-world Evidence-"

# Column G now holds the longest strings in the sheet, so size it to fit.
$ws.Columns.Item(7).EntireColumn.AutoFit()

# Leave the selection on the last-edited cell, as the author did.
$ws.Range("G12").Select()
